# edit.ps1 -- "suite chap 13 (PAS FINI FINALEMENT !)"
#
# Applies two changes to "Chapitre 13 - Vaelya.docx":
#
#   1. During the trial paragraph, Ser Barton no longer protests Valyra's
#      *innocence* -- he protests her *guilt* instead:
#         "...protesté de la prétendue innocence de Valyra."
#      becomes
#         "...protesté de la culpabilité de Valyra."
#
#   2. A brand new scene is appended after "...quelque chose ce soir.",
#      separated by a blank line, a centered "***" scene break, another
#      blank line, then a new narrative paragraph about the Sentinelles'
#      assignments, followed by a trailing blank paragraph. The `_GoBack`
#      bookmark (Word's "last edit location" marker) moves from the end of
#      the old final paragraph to the end of this new paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...la prétendue innocence de Valyra." -> "...la culpabilité de Valyra."
# ---------------------------------------------------------------------
$replaced = $d.Content.Find.Execute(
    "même s’il avait à grands cris protesté de la prétendue innocence de Valyra.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "même s’il avait à grands cris protesté de la culpabilité de Valyra.", 2)
if (-not $replaced) {
    throw "could not find the trial-verdict sentence to update"
}

# ---------------------------------------------------------------------
# 2) Append the new scene at the end of the chapter.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range

# The _GoBack bookmark currently marks the end of the chapter; remove it here,
# it is recreated at the end of the new paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -- blank paragraph --------------------------------------------------
$tail.InsertParagraphAfter()
$pBlank1 = $d.Paragraphs.Last
$pBlank1.Format.Alignment = 3

# -- centered "***" scene break ---------------------------------------
$pBlank1.Range.InsertParagraphAfter()
$pStar = $d.Paragraphs.Last
$pStar.Format.Alignment = 1
$pStar.Range.InsertAfter("***")

# -- blank paragraph ----------------------------------------------------
$pStar.Range.InsertParagraphAfter()
$pBlank2 = $d.Paragraphs.Last
$pBlank2.Format.Alignment = 3

# -- new narrative paragraph -------------------------------------------
$pBlank2.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Last
$pNew.Format.Alignment = 3

$pNew.Range.InsertAfter("L’élection n’eut pas lieu tout de suite. En revanche, un autre évènement arriva qui parvint presque à redonner le sourire à Vaelya, et lui fit oublier pendant quelques instants qu’elle avait échappé de peu à la pendaison : le Commandant allait annoncer les affectations des Sentinelles fraichement assermentées.")

$run2 = $pNew.Range.Duplicate
$run2.Collapse(0)
$run2.InsertAfter(" Vaelya avait discuté avec sa sœur afin d’avoir son avis. Valait-il mieux qu’elles restent ensemble, ou serait-il plus intéressant qu’elle devienne une Eclaireuse ?")

$run3 = $pNew.Range.Duplicate
$run3.Collapse(0)
$run3.InsertAfter(" Sa sœur l’avait convaincue de suivre son instinct et de demander à être affectée au corps des Eclaireurs. Cela ne les séparerait sûrement pas beaucoup, puisque les Sentinelles passaient de toute façon la majorité de leur temps à Grand-Roc.")

$run4 = $pNew.Range.Duplicate
$run4.Collapse(0)
$run4.InsertAfter(" Et Soldats et Eclaireurs partaient ")

$run5 = $pNew.Range.Duplicate
$run5.Collapse(0)
$run5.InsertAfter("souvent en ")

$run6 = $pNew.Range.Duplicate
$run6.Collapse(0)
$run6.InsertAfter("mission ensemble.")

$run7 = $pNew.Range.Duplicate
$run7.Collapse(0)
$run7.InsertAfter(" Valyra avait aussi ajouté pour embêter sa sœur que de toute façon, si elle avait l’intention d’être utile, le mieux était qu’elle choisisse le corps des Ouvriers, vu comment elle maniait les armes.")

# Recreate the _GoBack bookmark, collapsed, right at the end of the new
# paragraph's text (mirroring where Word leaves it after the last edit).
# A transient marker character works around collapsed-range bookmark
# placement right at the tail of the document.
$endMark = $pNew.Range.Duplicate
$endMark.Collapse(0)
$null = $endMark.MoveEnd(1, -1)
$endMark.InsertAfter([char]1)
$d.Bookmarks.Add("_GoBack", $endMark)
$endMark.Text = ""

# -- trailing blank paragraph -------------------------------------------
$pNew.Range.InsertParagraphAfter()
$pTrail = $d.Paragraphs.Last
$pTrail.Format.Alignment = 3

Write-Output "done"
